# Applies the cryptos-list price/volume refresh described in the commit diff.
# D-column "Price" values that are plain decimals (e.g. 0.7049) are written with a
# leading apostrophe so Excel stores them as literal text (quote-prefix), matching
# the original file where every Price/Volume cell is text, not a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.225.25'
$ws.Range('E2').Value = '  -0.55%  '

$ws.Range('D3').Value = '1.860.87'
$ws.Range('E3').Value = '  -1.25%  '

$ws.Range('D4').Value = '''1.002'

$ws.Range('D5').Value = '''0.7049'
$ws.Range('E5').Value = '  -1.06%  '

$ws.Range('D6').Value = '''242.22'
$ws.Range('E6').Value = '  -0.08%  '

$ws.Range('D7').Value = '''1.002'
$ws.Range('E7').Value = '  +0.04%  '

$ws.Range('D8').Value = '''0.07813'
$ws.Range('E8').Value = '  -2.75%  '

$ws.Range('D9').Value = '''0.3112'
$ws.Range('E9').Value = '  -0.50%  '

$ws.Range('D10').Value = '''24.24'
$ws.Range('E10').Value = '  -3.98%  '

$ws.Range('D11').Value = '''0.08003'
$ws.Range('E11').Value = '  -4.25%  '

$ws.Range('D12').Value = '1.868.90'
$ws.Range('E12').Value = '  -1.02%  '

$ws.Range('D13').Value = '''5.170'
$ws.Range('E13').Value = '  -1.44%  '

$ws.Range('D14').Value = '''93.37'
$ws.Range('E14').Value = '  +0.84%  '

$ws.Range('D15').Value = '''0.6944'
$ws.Range('E15').Value = '  -3.67%  '

$ws.Range('D16').Value = '''6.349'
$ws.Range('E16').Value = '  +0.77%  '

$ws.Range('D17').Value = '29.225.27'
$ws.Range('E17').Value = '  -0.62%  '

$ws.Range('D18').Value = '''0.000008281'
$ws.Range('E18').Value = '  -2.36%  '

$ws.Range('D19').Value = '''252.58'
$ws.Range('E19').Value = '  +4.70%  '

$ws.Range('D20').Value = '2.150.95'
$ws.Range('E20').Value = '  +0.33%  '

$ws.Range('E21').Value = '  -1.30%  '

$ws.Range('D22').Value = '''1.002'
$ws.Range('E22').Value = '  +0.07%  '

$ws.Range('D23').Value = '''7.519'
$ws.Range('E23').Value = '  -4.40%  '

$ws.Range('D25').Value = '''0.1552'
$ws.Range('E25').Value = '  -2.17%  '

$ws.Range('D26').Value = '''8.987'
$ws.Range('E26').Value = '  -0.99%  '

$ws.Range('D27').Value = '''159.35'
$ws.Range('E27').Value = '  -2.68%  '

$ws.Range('D28').Value = '''18.71'
$ws.Range('E28').Value = '  +0.71%  '

$ws.Range('E29').Value = '  -0.64%  '

$ws.Range('D30').Value = '''4.263'
$ws.Range('E30').Value = '  -1.79%  '

$ws.Range('D31').Value = '''4.271'
$ws.Range('E31').Value = '  -3.36%  '

$ws.Range('D32').Value = '''1.209'
$ws.Range('E32').Value = '  -0.26%  '

$ws.Range('E33').Value = '  -2.01%  '

$ws.Range('D34').Value = '''1.883'
$ws.Range('E34').Value = '  -3.56%  '

$ws.Range('D35').Value = '''0.7423'
$ws.Range('E35').Value = '  -0.94%  '

$ws.Range('D36').Value = '''1.154'
$ws.Range('E36').Value = '  -2.46%  '

$ws.Range('D37').Value = '''2.706'
$ws.Range('E37').Value = '  +0.22%  '

$ws.Range('E38').Value = '  -1.62%  '

$ws.Range('D39').Value = '1.247.77'
$ws.Range('E39').Value = '  -3.06%  '

$ws.Range('D40').Value = '''2.740'
$ws.Range('E40').Value = '  -0.28%  '

$ws.Range('D41').Value = '''6.311'
$ws.Range('E41').Value = '  -4.43%  '

$ws.Range('D42').Value = '''0.9039'
$ws.Range('E42').Value = '  +0.08%  '

$ws.Range('D43').Value = '''110.89'
$ws.Range('E43').Value = '  -0.60%  '

$ws.Range('D44').Value = '''71.55'
$ws.Range('E44').Value = '  -2.86%  '

$ws.Range('E45').Value = '  +0.00%  '

$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '''0.00000000128'
$ws.Range('E46').Value = '  -0.16%  '

$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').Value = '2.053.60'
$ws.Range('E47').Value = '  +1.16%  '

$ws.Range('D48').Value = '''0.5203'
$ws.Range('E48').Value = '  -0.30%  '

$ws.Range('D49').Value = '''1.778'
$ws.Range('E49').Value = '  -1.67%  '

$ws.Range('D50').Value = '''9.386'
$ws.Range('E50').Value = '  -1.31%  '

$ws.Range('D51').Value = '''1.010'
$ws.Range('E51').Value = '  +0.41%  '
